$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1084281.8
$ws.Range("I17").Value = 479.21054
$ws.Range("J17").Value = 1513286.9
$ws.Range("K17").Value = 1437.63162
$ws.Range("L17").Value = 4539860.699999999
$ws.Range("M17").Value = -1269.63162
$ws.Range("N17").Value = -4540196.699999999

$ws.Range("H32").Value = 1900.5
$ws.Range("I32").Value = 2001
$ws.Range("J32").Value = 1800
$ws.Range("K32").Value = 2001
$ws.Range("L32").Value = 1800
$ws.Range("M32").Value = -1675
$ws.Range("N32").Value = -2452

$ws.Range("H43").Value = 50009510
$ws.Range("I43").Value = 1767
$ws.Range("J43").Value = 71441400
$ws.Range("K43").Value = 1767
$ws.Range("L43").Value = 71441400
$ws.Range("M43").Value = -1698
$ws.Range("N43").Value = -71441538

$ws.Range("H70").Value = 1550
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1550
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 4650
$ws.Range("M70").Value = ""
$ws.Range("N70").Value = -5190

$ws.Range("H73").Value = 1550
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1550
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 4650
$ws.Range("M73").Value = ""
$ws.Range("N73").Value = -6522

$ws.Range("H76").Value = 3085.7144
$ws.Range("I76").Value = 3100
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 3100
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -2785
$ws.Range("N76").Value = -3630

$ws.Range("H79").Value = 3085.7144
$ws.Range("I79").Value = 3100
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 3100
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -2008
$ws.Range("N79").Value = -5184

$ws.Range("H112").Value = 5726
$ws.Range("I112").Value = 350
$ws.Range("J112").Value = 6467.517
$ws.Range("K112").Value = 1050
$ws.Range("L112").Value = 19402.551
$ws.Range("M112").Value = 58
$ws.Range("N112").Value = -21618.551

$ws.Range("H129").Value = 1164.2295
$ws.Range("I129").Value = 472.64706
$ws.Range("J129").Value = 1431.4318
$ws.Range("K129").Value = 1417.94118
$ws.Range("L129").Value = 4294.2954
$ws.Range("M129").Value = 3582.05882
$ws.Range("N129").Value = -14294.2954

$ws.Range("H135").Value = 2343.7368
$ws.Range("I135").Value = 968.5833
$ws.Range("J135").Value = 4701.143
$ws.Range("K135").Value = 8717.2497
$ws.Range("L135").Value = 42310.287
$ws.Range("M135").Value = -6182.2497
$ws.Range("N135").Value = -47380.287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1261.7222
$ws.Range("I74").Value = 1037.4117
$ws.Range("J74").Value = 1462.421
$ws.Range("K74").Value = 1037.4117
$ws.Range("L74").Value = 1462.421
$ws.Range("M74").Value = -163.4117000000001
$ws.Range("N74").Value = -3210.421

$ws.Range("H77").Value = 1261.7222
$ws.Range("I77").Value = 1037.4117
$ws.Range("J77").Value = 1462.421
$ws.Range("K77").Value = 5187.058500000001
$ws.Range("L77").Value = 7312.105
$ws.Range("M77").Value = -819.058500000001
$ws.Range("N77").Value = -16048.105

$ws.Range("H110").Value = 77484.5
$ws.Range("I110").Value = 102411.445
$ws.Range("J110").Value = 2703.6667
$ws.Range("K110").Value = 102411.445
$ws.Range("L110").Value = 2703.6667
$ws.Range("M110").Value = -100366.445
$ws.Range("N110").Value = -6793.6667

$ws.Range("H132").Value = 2453.8833
$ws.Range("I132").Value = 1407.1333
$ws.Range("J132").Value = 5594.1333
$ws.Range("K132").Value = 4221.3999
$ws.Range("L132").Value = 16782.3999
$ws.Range("M132").Value = -1691.3999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 25002864
$ws.Range("I105").Value = 25002864
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 25002864
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -25001117
$ws.Range("N105").Value = ""

$ws.Range("H134").Value = 2762.7727
$ws.Range("I134").Value = 2413.1155
$ws.Range("J134").Value = 3267.8333
$ws.Range("K134").Value = 7239.3465
$ws.Range("L134").Value = 9803.499899999999
$ws.Range("M134").Value = -4704.3465
$ws.Range("N134").Value = -14873.4999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 1826.2
$ws.Range("I33").Value = 1826.2
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 1826.2
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -1447.2
$ws.Range("N33").Value = ""

$ws.Range("H117").Value = 44027
$ws.Range("I117").Value = 19950
$ws.Range("J117").Value = 49377.445
$ws.Range("K117").Value = 19950
$ws.Range("L117").Value = 49377.445
$ws.Range("M117").Value = -15361
$ws.Range("N117").Value = -58555.445

$ws.Range("H132").Value = 40231690
$ws.Range("I132").Value = 45456056
$ws.Range("J132").Value = 23812240
$ws.Range("K132").Value = 136368168
$ws.Range("L132").Value = 71436720
$ws.Range("M132").Value = -136365638
$ws.Range("N132").Value = -71441780

$ws.Range("H139").Value = 29318.092
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 29318.092
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 29318.092
$ws.Range("N139").Value = -39598.092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1741.375
$ws.Range("I5").Value = 588
$ws.Range("J5").Value = 2318.0625
$ws.Range("K5").Value = 1764
$ws.Range("L5").Value = 6954.1875
$ws.Range("M5").Value = -1652
$ws.Range("N5").Value = -7178.1875

$ws.Range("H14").Value = 166666850
$ws.Range("I14").Value = 166666850
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 500000550
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -500000377

$ws.Range("H135").Value = 1741.375
$ws.Range("I135").Value = 588
$ws.Range("J135").Value = 2318.0625
$ws.Range("K135").Value = 5292
$ws.Range("L135").Value = 20862.5625
$ws.Range("M135").Value = -2757
$ws.Range("N135").Value = -25932.5625

$ws.Range("H138").Value = 3834.15
$ws.Range("I138").Value = 1459.75
$ws.Range("J138").Value = 7395.75
$ws.Range("K138").Value = 4379.25
$ws.Range("L138").Value = 22187.25
$ws.Range("M138").Value = 760.75
$ws.Range("N138").Value = -32467.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4000
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 5000
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = -888
$ws.Range("N5").Value = -5224

$ws.Range("H113").Value = 2088.1304
$ws.Range("I113").Value = 2020.7778
$ws.Range("J113").Value = 2330.6
$ws.Range("K113").Value = 2020.7778
$ws.Range("L113").Value = 2330.6
$ws.Range("M113").Value = 149.2221999999999
$ws.Range("N113").Value = -6670.6

$ws.Range("H122").Value = 1585.6666
$ws.Range("I122").Value = 1128.5
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 3385.5
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -935.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 6843.5264
$ws.Range("I93").Value = 8452.643
$ws.Range("J93").Value = 2338
$ws.Range("K93").Value = 8452.643
$ws.Range("L93").Value = 2338
$ws.Range("M93").Value = -7204.643
$ws.Range("N93").Value = -4834

$ws.Range("H132").Value = 3015.1
$ws.Range("I132").Value = 3090.7334
$ws.Range("J132").Value = 2939.4666
$ws.Range("K132").Value = 9272.200199999999
$ws.Range("L132").Value = 8818.399800000001
$ws.Range("M132").Value = -6742.200199999999
$ws.Range("N132").Value = -13878.3998

$ws.Range("H136").Value = 27781180
$ws.Range("I136").Value = 5702
$ws.Range("J136").Value = 41668920
$ws.Range("K136").Value = 17106
$ws.Range("L136").Value = 125006760
$ws.Range("M136").Value = -14556
$ws.Range("N136").Value = -125011860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 500000500
$ws.Range("I32").Value = 1026
$ws.Range("J32").Value = 1000000000
$ws.Range("K32").Value = 1026
$ws.Range("L32").Value = 1000000000
$ws.Range("M32").Value = -709
$ws.Range("N32").Value = -1000000634

$ws.Range("H100").Value = 1063.8422
$ws.Range("I100").Value = 1385.5555
$ws.Range("J100").Value = 774.3
$ws.Range("K100").Value = 2771.111
$ws.Range("L100").Value = 1548.6
$ws.Range("M100").Value = -2230.111
$ws.Range("N100").Value = -2630.6

$ws.Range("H111").Value = 42879.668
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 42879.668
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 42879.668
$ws.Range("N111").Value = -51059.668

$ws.Range("H113").Value = 1345.826
$ws.Range("I113").Value = 1530.2667
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 4590.800099999999
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -2420.800099999999
$ws.Range("N113").Value = -7340

$ws.Range("H122").Value = 1802.0851
$ws.Range("I122").Value = 1763.5555
$ws.Range("J122").Value = 1928.1818
$ws.Range("K122").Value = 5290.666499999999
$ws.Range("L122").Value = 5784.5454
$ws.Range("M122").Value = -2840.666499999999
$ws.Range("N122").Value = -10684.5454

$ws.Range("H132").Value = 5211279.5
$ws.Range("I132").Value = 4576.9165
$ws.Range("J132").Value = 8335301
$ws.Range("K132").Value = 13730.7495
$ws.Range("L132").Value = 25005903
$ws.Range("M132").Value = -11200.7495
$ws.Range("N132").Value = -25010963

$ws.Range("H136").Value = 3414.1428
$ws.Range("I136").Value = 3124.4783
$ws.Range("J136").Value = 3969.3333
$ws.Range("K136").Value = 9373.4349
$ws.Range("L136").Value = 11907.9999
$ws.Range("M136").Value = -6823.4349
$ws.Range("N136").Value = -17007.9999
